$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.682.86"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "3.794.36"
$ws.Range("E3").Value = "  +1.59%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'595.80"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "'167.16"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").Value = "3.792.01"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").Value = "'6.32"
$ws.Range("E11").Value = "  -1.24%  "

$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("D14").Value = "'36.09"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "4.434.25"
$ws.Range("E15").Value = "  +1.50%  "

$ws.Range("D16").Value = "3.807.17"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.654.80"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.41"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'9.99"
$ws.Range("E21").Value = "  -6.02%  "

$ws.Range("D22").Value = "'458.56"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "'0.697"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "'0.0000155"
$ws.Range("E24").Value = "  +6.27%  "

$ws.Range("D25").Value = "'83.37"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").Value = "'12.03"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").Value = "'2.11"
$ws.Range("E27").Value = "  -2.45%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'10.00"

$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("E31").Value = "  +4.43%  "

$ws.Range("D32").Value = "'7.23"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "'29.58"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").Value = "'9.08"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").Value = "3.734.42"
$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("D37").Value = "'0.0999"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").Value = "'3.38"
$ws.Range("E38").Value = "  -1.15%  "

$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").Value = "'0.990"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").Value = "'5.76"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'44.49"
$ws.Range("E44").Value = "  +1.17%  "

$ws.Range("D45").Value = "'48.03"
$ws.Range("E45").Value = "  +2.77%  "

$ws.Range("D46").Value = "'0.298"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'149.07"
$ws.Range("E47").Value = "  +3.54%  "

$ws.Range("D48").Value = "'8.27"
$ws.Range("E48").Value = "  -1.46%  "

$ws.Range("D49").Value = "'393.47"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.82"
$ws.Range("E50").Value = "  -4.03%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'26.57"
$ws.Range("E51").Value = "  +6.46%  "
